$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.247.76'
$ws.Range('E2').Value = '  -0.13%  '
$ws.Range('D3').Value = '3.533.99'
$ws.Range('E3').Value = '  +3.23%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Formula = '''595.51'
$ws.Range('E5').Value = '  +0.22%  '
$ws.Range('D6').Formula = '''139.05'
$ws.Range('E6').Value = '  +3.23%  '
$ws.Range('D7').Value = '3.534.09'
$ws.Range('E7').Value = '  +3.31%  '
$ws.Range('E8').Value = '  +0.23%  '
$ws.Range('D9').Formula = '''0.495'
$ws.Range('E9').Value = '  +1.19%  '
$ws.Range('E10').Value = '  +3.73%  '
$ws.Range('D11').Formula = '''7.18'
$ws.Range('E11').Value = '  -2.97%  '
$ws.Range('D12').Formula = '''0.389'
$ws.Range('E12').Value = '  +3.76%  '
$ws.Range('D13').Value = '4.134.13'
$ws.Range('E13').Value = '  +3.35%  '
$ws.Range('E14').Value = '  +4.32%  '
$ws.Range('D15').Formula = '''26.95'
$ws.Range('E15').Value = '  +2.36%  '
$ws.Range('D16').Value = '3.522.02'
$ws.Range('E16').Value = '  +2.60%  '
$ws.Range('E17').Value = '  +1.39%  '
$ws.Range('D18').Value = '65.147.00'
$ws.Range('E18').Value = '  -0.18%  '
$ws.Range('D19').Formula = '''10.19'
$ws.Range('E19').Value = '  +2.39%  '
$ws.Range('D20').Formula = '''5.84'
$ws.Range('E20').Value = '  +2.34%  '
$ws.Range('D21').Formula = '''14.17'
$ws.Range('E21').Value = '  +4.17%  '
$ws.Range('D22').Formula = '''394.79'
$ws.Range('E22').Value = '  +1.38%  '
$ws.Range('D23').Formula = '''0.569'
$ws.Range('E23').Value = '  +5.20%  '
$ws.Range('D24').Formula = '''74.43'
$ws.Range('E24').Value = '  +1.91%  '
$ws.Range('D25').Value = '3.679.04'
$ws.Range('E25').Value = '  +3.17%  '
$ws.Range('D26').Formula = '''0.999'
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('D27').Formula = '''0.0000113'
$ws.Range('E27').Value = '  +8.84%  '
$ws.Range('D28').Formula = '''7.77'
$ws.Range('E28').Value = '  +9.10%  '
$ws.Range('D29').Formula = '''1.00'
$ws.Range('E29').Value = '  +0.24%  '
$ws.Range('E30').Value = '  +0.83%  '
$ws.Range('D31').Formula = '''8.24'
$ws.Range('E31').Value = '  +1.50%  '
$ws.Range('D32').Value = '3.554.33'
$ws.Range('E32').Value = '  +3.72%  '
$ws.Range('E33').Value = '  +0.00%  '
$ws.Range('D34').Formula = '''23.78'
$ws.Range('E34').Value = '  +5.60%  '
$ws.Range('E35').Value = '  +1.09%  '
$ws.Range('D36').Formula = '''1.24'
$ws.Range('E36').Value = '  +0.47%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').Formula = '''170.42'
$ws.Range('E37').Value = '  -1.34%  '
$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D38').Formula = '''6.97'
$ws.Range('E38').Value = '  +2.55%  '
$ws.Range('E39').Value = '  +1.45%  '
$ws.Range('D40').Formula = '''4.89'
$ws.Range('E40').Value = '  +2.03%  '
$ws.Range('D41').Formula = '''0.0798'
$ws.Range('E41').Value = '  +4.33%  '
$ws.Range('B42').Value = 'Mantle'
$ws.Range('C42').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D42').Formula = '''0.821'
$ws.Range('E42').Value = '  +1.33%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').Formula = '''26.61'
$ws.Range('E43').Value = '  +22.46%  '
$ws.Range('D44').Formula = '''42.66'
$ws.Range('E44').Value = '  -2.15%  '
$ws.Range('E45').Value = '  -0.14%  '
$ws.Range('D46').Formula = '''4.41'
$ws.Range('E46').Value = '  +1.15%  '
$ws.Range('E47').Value = '  +9.46%  '
$ws.Range('D48').Formula = '''1.66'
$ws.Range('E48').Value = '  +3.59%  '
$ws.Range('D49').Formula = '''6.80'
$ws.Range('E49').Value = '  +4.38%  '
$ws.Range('D50').Value = '2.344.73'
$ws.Range('E50').Value = '  +6.78%  '
$ws.Range('D51').Formula = '''2.12'
$ws.Range('E51').Value = '  +0.04%  '
